$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete old data rows (rows 2-17), keep header row 1
$ws.Range("A2:C17").ClearContents()

# New data values (rows 2-4)
$data = @(
    @(4,1,4),
    @(4,5,3),
    @(4,5,5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

$ws.Range("C4").Select()
